# CDS Study filter fixes
# - Renames the "Cases" row into a "Participants" row (keeps the participant
#   query, which had already been living in that row's query column).
# - Fixes the broken/duplicated StatQuery text used in column C.
# - Bumps the base font size used across the grid from 12 to 15 and drops the
#   wrap-text default style so column A/D/E read cleanly, while column B/C
#   (the long query text) keep word-wrap.
# - Re-sizes columns/rows to fit the new font and re-selects E9.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Query / label text blocks (kept verbatim from the workbook's query bank)
# ---------------------------------------------------------------------------

$participantQuery = @"
MATCH (s:study)<--(p:participant)
WHERE s.study_name in ["University of Texas PDX Development and Trial Center Grant"]
OPTIONAL MATCH (p)<--(samp:sample)
WITH p, s, collect(distinct samp.sample_id) as samp
RETURN   
 coalesce(p.participant_id,'') as ``Participant ID``,
 coalesce(s.study_name, '') as ``Study Name``,
 coalesce(s.phs_accession,'') as ``Accession``,
 coalesce(p.gender,'') as ``Gender``,
 coalesce(apoc.text.join(samp, ','), '') as ``Samples``
 ORDER By p.participant_id LIMIT 100
"@

$sampleQuery = @"
MATCH (s:study)<--(p:participant)<--(samp:sample)
WHERE s.study_name in ["University of Texas PDX Development and Trial Center Grant"]
WITH p,s,samp,COLLECT(DISTINCT samp.sample_tumor_status) as tumor
RETURN  
 coalesce(samp.sample_id, '') as ``Sample ID``,
 coalesce(p.participant_id,'') as ``Participant ID``,
 coalesce(s.study_name, '') as ``Study Name``,
 coalesce(s.phs_accession,'') as ``Accession``,
coalesce(samp.sample_tumor_status,'') as ``Tumor``,
coalesce(samp.sample_type,'') as ``Analyte Type``
ORDER By samp.sample_id LIMIT 100
"@

$fileQuery = @"
MATCH (s:study)<--(p:participant)
WHERE s.study_name in ["University of Texas PDX Development and Trial Center Grant"]
OPTIONAL MATCH (p)<--(samp:sample)
OPTIONAL MATCH (p)<--(diag:diagnosis)
OPTIONAL MATCH (samp)<--(f:file)
WITH DISTINCT p,s,samp,f,diag
RETURN 
    coalesce(f.file_name, '') as ``File Name``,
    coalesce(s.study_name, '') as ``Study Name``,
    coalesce(s.phs_accession,'') as ``Accession``,
    coalesce(p.participant_id,'') as ``Participant ID``,
    coalesce(samp.sample_id, '') as ``Sample ID``,
    coalesce(f.file_type, '') as ``File Type``
   ORDER By f.file_name LIMIT 100
"@

# The dbExcel "StatQuery" column used the same (broken / duplicated) Cypher
# snippet for every tab. Fix it up -- drop the duplicated middle section --
# and keep it as the shared value for all three data rows.
$statQuery = @"
MATCH (s:study)<--(p:participant)
OPTIONAL MATCH (p)<--(samp:sample)MATCH (s:study)<--(p:participant)
OPTIONAL MATCH (p)<--(samp:sample)
OPTIONAL MATCH (p)<--(diag:diagnosis)
OPTIONAL MATCH (samp)<--(f:file)
WITH DISTINCT samp,diag,s,p,f
WHERE s.study_name in ["University of Texas PDX Development and Trial Center Grant"]
RETURN
    count(distinct s) AS Studies,
    count(distinct p) AS Participants,
    count(distinct samp) AS Samples,
    count(distinct f) AS ``Files``   
"@

$neo4jFile = "TC11_CDS_Filter_Study-UniversityofTexas_Neo4jData.xlsx"
$webFile = "TC11_CDS_Filter_Study-UniversityofTexas_WebData.xlsx"

# ---------------------------------------------------------------------------
# Grid values
# ---------------------------------------------------------------------------

$ws.Range("A1").Value = "TabName"
$ws.Range("B1").Value = "query"
$ws.Range("C1").Value = "StatQuery"
$ws.Range("D1").Value = "dbExcel"
$ws.Range("E1").Value = "WebExcel"

# Row 2: CasesTab -> ParticipantsTab (query column already held the
# participant query, so the row is simply relabelled & fixed up)
$ws.Range("A2").Value = "ParticipantsTab"
$ws.Range("B2").Value = $participantQuery
$ws.Range("C2").Value = $statQuery
$ws.Range("D2").Value = $neo4jFile
$ws.Range("E2").Value = $webFile

# Row 3: SamplesTab
$ws.Range("A3").Value = "SamplesTab"
$ws.Range("B3").Value = $sampleQuery
$ws.Range("C3").Value = $statQuery
$ws.Range("D3").Value = $neo4jFile
$ws.Range("E3").Value = $webFile

# Row 4: FilesTab
$ws.Range("A4").Value = "FilesTab"
$ws.Range("B4").Value = $fileQuery
$ws.Range("C4").Value = $statQuery
$ws.Range("D4").Value = $neo4jFile
$ws.Range("E4").Value = $webFile

# ---------------------------------------------------------------------------
# Row heights (unchanged for rows 2-3, row 4 shrinks now its text is tidier)
# ---------------------------------------------------------------------------

$ws.Rows.Item(2).RowHeight = 213.75
$ws.Rows.Item(3).RowHeight = 209.25
$ws.Rows.Item(4).RowHeight = 222.75

# ---------------------------------------------------------------------------
# Fonts / wrap text
# ---------------------------------------------------------------------------

# Whole used grid goes to size 15 ...
$ws.Range("A1:E6").Font.Size = 15

# ... but only the long query columns (B/C) keep word-wrap.
$ws.Range("B2:C4").WrapText = $true
$ws.Range("B5:C6").WrapText = $true

# Columns A, D, E (and the header row) should NOT wrap.
$ws.Range("A1:E1").WrapText = $false
$ws.Range("A2:A4").WrapText = $false
$ws.Range("D2:E4").WrapText = $false

# ---------------------------------------------------------------------------
# Column widths
# ---------------------------------------------------------------------------

$ws.Columns.Item(1).ColumnWidth = 21.29
$ws.Columns.Item(2).ColumnWidth = 75.71
$ws.Columns.Item(3).ColumnWidth = 74.86
$ws.Columns.Item(4).ColumnWidth = 52
$ws.Columns.Item(5).ColumnWidth = 73.14

# ---------------------------------------------------------------------------
# Selection (last thing the author clicked before saving)
# ---------------------------------------------------------------------------

$ws.Range("E9").Select()
